$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formations_IED")

# Row 7: "Responsable Achats Durables" -> add Programme_Secondaire (C) and Modules_Clés (D)
$ws.Range("C7").Value = "Mastère Spécialisé Supply Chain Durable ou Certification RSE"
$ws.Range("D7").Value = "Durabilité & Achats et Supply Chain, Cartographie des risques fournisseurs, Économie circulaire, Analyse cycle de vie, Normes ISO 20400"

# Row 8: "Chargé de Mission Biodiversité" -> add Programme_Secondaire (C) and Modules_Clés (D)
$ws.Range("C8").Value = "Master Écologie & Biodiversité ou MSc Conservation"
$ws.Range("D8").Value = "Préservation des écosystèmes, Mesure de la performance environnementale, Reporting biodiversité, Comptabilités multi-capitaux, Valorisation des externalités"

# Match formatting of the other data rows (style used by C2:D6) for the newly filled cells
$ws.Range("C2:D2").Copy()
$ws.Range("C7:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
